# Updated symbol list on Mon Dec 19 21:19:22 UTC 2022 with GitHub Actions
# Refresh cryptocurrency ranking data (price/volume columns) and re-rank a
# handful of coins whose 24h order shifted.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking values stored as text in the
# source workbook. A leading apostrophe forces Excel to keep them as text
# (matching the original inlineStr cells) instead of auto-converting to
# numbers.
$ws.Range("D2").Formula = "'243.63"
$ws.Range("D3").Formula = "'21.50"
$ws.Range("D5").Formula = "'0.05602"
$ws.Range("D6").Formula = "'3.363"
$ws.Range("D7").Formula = "'6.369"
$ws.Range("D8").Formula = "'0.8055"
$ws.Range("D9").Formula = "'0.9336"
$ws.Range("D10").Formula = "'0.1424"
$ws.Range("D11").Formula = "'0.07290"
$ws.Range("D12").Formula = "'0.03113"
$ws.Range("D13").Formula = "'0.03053"
$ws.Range("D14").Formula = "'0.09276"
$ws.Range("D15").Formula = "'3.570"
$ws.Range("D16").Formula = "'0.001639"
$ws.Range("D17").Formula = "'0.04696"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Formula = "'0.006404"
$ws.Range("E18").Value = "17TigerCashTCH"
$ws.Range("B19").Value = "HotbitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D19").Formula = "'0.004989"
$ws.Range("E19").Value = "18HotbitTokenHTB"
$ws.Range("B20").Value = "BitKan"
$ws.Range("C20").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D20").Formula = "'0.001042"
$ws.Range("E20").Value = "19BitKanKAN"
$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D21").Formula = "'0.0001499"
$ws.Range("E21").Value = "20NitroExNTX"
$ws.Range("B22").Value = "UpBots"
$ws.Range("C22").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D22").Formula = "'0.0003099"
$ws.Range("E22").Value = "21UpBotsUBXT"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").Formula = "'3.757"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").Formula = "'2.094"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "One"
$ws.Range("C25").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D25").Formula = "'0.01120"
$ws.Range("E25").Value = "24OneONEBestin24h"
$ws.Range("D26").Formula = "'0.3260"
$ws.Range("D27").Formula = "'0.1281"
$ws.Range("D40").Formula = "'0.03919"
$ws.Range("D41").Formula = "'0.006884"
$ws.Range("D43").Formula = "'0.1033"
$ws.Range("D44").Formula = "'0.008491"
$ws.Range("D45").Formula = "'0.00005927"
$ws.Range("D47").Formula = "'0.0005498"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
$ws.Range("D48").Formula = "'0.6822"
$ws.Range("D49").Formula = "'0.08439"
$ws.Range("E49").Value = "48BOLOBOLO"
